$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 415, shifting existing rows 415-513 down to 416-514.
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row 415 with the new data.
$ws.Cells.Item(415, 1).Value = 3
$ws.Cells.Item(415, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(415, 3).Value = "Coquimbo"
$ws.Cells.Item(415, 4).Value = 44782
$ws.Cells.Item(415, 5).Value = 5
$ws.Cells.Item(415, 6).Value = 100112021
$ws.Cells.Item(415, 7).Value = "Ají"
$ws.Cells.Item(415, 8).Value = "Inferno"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 73
$ws.Cells.Item(415, 11).Value = 13000
$ws.Cells.Item(415, 12).Value = 14000
$ws.Cells.Item(415, 13).Value = 13521
$ws.Cells.Item(415, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(415, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(415, 16).Value = 901
$ws.Cells.Item(415, 17).Value = 15
$ws.Cells.Item(415, 18).Value = "Hortaliza"
